# Fix Navbar highlight tab & add codewords<->students mapping (col A = codeword
# hyperlink "mailto:" addresses, col B = student names). Two trailing codeword
# rows (6 & 7) carry the Hyperlink cell style but no data/link yet, since there
# aren't enough students to match those codewords.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (codewords) / Column B (student names) ------------------------
$ws.Range("A1").Value = "S530742@nwmissouri.edu"
$ws.Range("B1").Value = "Naveen"

$ws.Range("A2").Value = "S530742@nwmissouri.edu"
$ws.Range("B2").Value = "Naveen"

$ws.Range("A3").Value = "S530742@nwmissouri.edu"
$ws.Range("B3").Value = "Naveen"

$ws.Range("A4").Value = "S530742@nwmissouri.edu"
$ws.Range("B4").Value = "Naveen"

$ws.Range("A5").Value = "bob@bob.com"
$ws.Range("B5").Value = "Bobby"

# --- Rows 6 & 7: reserve the Hyperlink look (font/underline/color) for extra
# codewords that don't have a matching student yet, without linking them to
# anything. Add throwaway hyperlinks first (so the engine mints the
# Hyperlink cell style), then strip the links back off before the real
# hyperlinks go on A1:A5 -- Range.Hyperlinks.Delete() clears every hyperlink
# on the sheet, so this has to happen before A1:A5 get their real links.
$ws.Range("A6").Value = "placeholder"
$ws.Range("A7").Value = "placeholder"
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:placeholder@placeholder.com")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:placeholder@placeholder.com")
$ws.Range("A6").Hyperlinks.Delete()
$ws.Range("A6").Value = ""
$ws.Range("A7").Value = ""

# --- Real hyperlinks for the codewords that do have a student -------------
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:S530742@nwmissouri.edu")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:S530742@nwmissouri.edu")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:S530742@nwmissouri.edu")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:S530742@nwmissouri.edu")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:bob@bob.com")

# --- Selection / active cell matches the tab the author ended on ----------
$ws.Range("B5").Select()
